# Add new rows of ORA Error data (rows 20-23) to the Online worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: Date (as date serial values) and Error Count
$newData = @(
    @{ Row = 20; Date = 45982; Count = 74 },
    @{ Row = 21; Date = 45987; Count = 70 },
    @{ Row = 22; Date = 45985; Count = 85 },
    @{ Row = 23; Date = 45986; Count = 103 }
)

foreach ($entry in $newData) {
    $r = $entry.Row
    # Column A: date value, formatted like the existing date column (copy style from A19)
    $ws.Cells.Item($r, 1).Value = $entry.Date
    $ws.Cells.Item(19, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Cells.Item($r, 1).Value = $entry.Date

    # Column B: error count
    $ws.Cells.Item($r, 2).Value = $entry.Count
}
$excel.CutCopyMode = $false

# Update the selection / view state to match the edited document
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("H25").Select()
